$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the two retired question rows (Franco/Spain, and the
#    "president of the USA"/Joe Biden row), via native row deletes so the
#    remaining shared strings keep their relative order.
# ---------------------------------------------------------------------------

# "General Franco became leader..." / "Spain" / "Location" -> row 2
$ws.Rows(2).Delete()

# "Who is the president of the USA?" / "Joe Biden" / "Person" -> now row 5
$ws.Rows(5).Delete()

# Sheet now reads:
#   1 Question | Answer | Category
#   2 Innsbruck... | Austria | Location
#   3 Euro... | 1999 | Year
#   4 Schumacher (first) | 1994 | Year
#   5 F1 2022 | Max Verstappen | Person

# ---------------------------------------------------------------------------
# 2) Reorder: move the "F1 2022 champion" row up above the Euro question,
#    and push the Euro question down below the Schumacher "first title" row.
# ---------------------------------------------------------------------------

# Move row 5 (F1 2022 / Max Verstappen / Person) up to row 3.
$ws.Rows(3).Insert()
$ws.Range("A6:C6").Cut($ws.Range("A3:C3"))
$ws.Rows(6).Delete()

# Sheet now reads:
#   1 Question | Answer | Category
#   2 Innsbruck... | Austria | Location
#   3 F1 2022 | Max Verstappen | Person
#   4 Euro... | 1999 | Year
#   5 Schumacher (first) | 1994 | Year

# Move the Euro question (row 4) down to after the Schumacher row (row 5).
$ws.Rows(6).Insert()
$ws.Range("A4:C4").Cut($ws.Range("A6:C6"))
$ws.Rows(4).Delete()

# Sheet now reads:
#   1 Question | Answer | Category
#   2 Innsbruck... | Austria | Location
#   3 F1 2022 | Max Verstappen | Person
#   4 Schumacher (first) | 1994 | Year
#   5 Euro... | 1999 | Year

# ---------------------------------------------------------------------------
# 3) Insert the new "2nd title" question between the first Schumacher
#    question and the Euro question, then append the rest of the new
#    Schumacher/World Cup/Chelsea questions at the bottom of the sheet.
# ---------------------------------------------------------------------------

$ws.Rows(5).Insert()
$ws.Cells.Item(5, 1).Value = "When did Miachel Schumacher win his 2nd F1 World Drivers Title?"
$ws.Cells.Item(5, 2).Value = 1995
$ws.Cells.Item(5, 3).Value = "Year"

$ws.Cells.Item(7, 1).Value = "When did Miachel Schumacher win his 3rd F1 World Drivers Title?"
$ws.Cells.Item(7, 2).Value = 2000
$ws.Cells.Item(7, 3).Value = "Year"

$ws.Cells.Item(8, 1).Value = "When did Miachel Schumacher win his 4th F1 World Drivers Title?"
$ws.Cells.Item(8, 2).Value = 2001
$ws.Cells.Item(8, 3).Value = "Year"

$ws.Cells.Item(9, 1).Value = "When did Miachel Schumacher win his 5th F1 World Drivers Title?"
$ws.Cells.Item(9, 2).Value = 2002
$ws.Cells.Item(9, 3).Value = "Year"

$ws.Cells.Item(10, 1).Value = "When did Miachel Schumacher win his 6th F1 World Drivers Title?"
$ws.Cells.Item(10, 2).Value = 2003
$ws.Cells.Item(10, 3).Value = "Year"

$ws.Cells.Item(11, 1).Value = "When did Miachel Schumacher win his 7th F1 World Drivers Title?"
$ws.Cells.Item(11, 2).Value = 2004
$ws.Cells.Item(11, 3).Value = "Year"

# NOTE: the "Chelsea" question text was introduced (first typed) before the
# "World Cup" question text, even though it ends up one row further down the
# sheet -- write the values in that order so new shared-string entries land
# in the same sequence as the authored workbook.
$ws.Cells.Item(13, 1).Value = "When has Chelsea last won the Champions League?"
$ws.Cells.Item(12, 1).Value = "When was the World Cup last won by Italy?"

$ws.Cells.Item(12, 2).Value = 2006
$ws.Cells.Item(12, 3).Value = "Year"

$ws.Cells.Item(13, 2).Value = 2021
$ws.Cells.Item(13, 3).Value = "Year"

# ---------------------------------------------------------------------------
# 4) Match the author's final selection/active cell.
# ---------------------------------------------------------------------------
$ws.Range("G17").Select() | Out-Null
